$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.107.53'
$ws.Range('E2').Value = '  -2.82%  '
$ws.Range('D3').Value = '1.651.86'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.14'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4788'
$ws.Range('E7').Value = '  -8.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2613'
$ws.Range('E8').Value = '  -4.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.05961'
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07070'
$ws.Range('E10').Value = '  -1.43%  '
$ws.Range('D11').Value = '1.660.85'
$ws.Range('E11').Value = '  -4.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.36'
$ws.Range('E12').Value = '  -4.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6166'
$ws.Range('E13').Value = '  -3.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.562'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '72.81'
$ws.Range('E15').Value = '  -5.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '25.125.99'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.34'
$ws.Range('E19').Value = '  -3.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006530'
$ws.Range('E20').Value = '  -3.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.408'
$ws.Range('E21').Value = '  +3.29%  '
$ws.Range('D22').Value = '1.872.90'
$ws.Range('E22').Value = '  -4.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.439'
$ws.Range('E23').Value = '  -2.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.259'
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '132.58'
$ws.Range('E25').Value = '  -4.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.66'
$ws.Range('E26').Value = '  -3.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.384'
$ws.Range('E27').Value = '  -8.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.703'
$ws.Range('E28').Value = '  -3.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '101.91'
$ws.Range('E29').Value = '  -3.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.807'
$ws.Range('E30').Value = '  -3.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07875'
$ws.Range('E31').Value = '  -4.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.517'
$ws.Range('E32').Value = '  -4.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04588'
$ws.Range('E33').Value = '  -0.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.609'
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9398'
$ws.Range('E35').Value = '  -4.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.5824'
$ws.Range('E36').Value = '  -5.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.618'
$ws.Range('E37').Value = '  -2.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01537'
$ws.Range('E38').Value = '  -4.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.001'
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8360'
$ws.Range('E40').Value = '  +12.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.831'
$ws.Range('E41').Value = '  -5.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.81'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3680'
$ws.Range('E43').Value = '  -4.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.832'
$ws.Range('E44').Value = '  -3.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1124'
$ws.Range('E45').Value = '  -0.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.025'
$ws.Range('E46').Value = '  -3.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05145'
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '52.04'
$ws.Range('E48').Value = '  -4.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.44'
$ws.Range('E49').Value = '  -3.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.001'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.288'
$ws.Range('E51').Value = '  -3.70%  '
